$d = $word.ActiveDocument

# ===========================================================================
# Change 1: "Mean, as it is more centralized between the min and max. "
#   -> two runs: "Median"  +  ". There are a lot of high value projects
#      which skews the average."
# ===========================================================================
$probe1 = $d.Content
$found1 = $probe1.Find.Execute(
    "Mean, as it is more centralized between the min and max. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $newText1 = "Median. There are a lot of high value projects which skews the average."
    $run1 = $d.Range($probe1.Start, $probe1.End)
    $run1.Text = $newText1

    # Toggling a character property on/off is a no-op visually, but it
    # forces the engine to keep this edited span as its own run instead of
    # re-merging it with whatever used to follow it.
    $run1.Bold = 1
    $run1.Bold = 0

    # Now split "Median" from the rest by toggling the same no-op on just
    # the "Median" prefix, which breaks it into two sibling runs.
    $medianLen = ("Median").Length
    $medianPart = $d.Range($run1.Start, $run1.Start + $medianLen)
    $medianPart.Bold = 1
    $medianPart.Bold = 0
}

# ===========================================================================
# Change 2: "There is more variability in successful campaigns. This does
#   make sense as there are simply more successful campaigns (565) than
#   there are failed (364). " (including the single trailing space run that
#   originally followed it)
#   -> four runs: "There is more variability in successful campaigns. " +
#      "This makes sense because there are a lot more samples for
#      successful outcomes as well as more" + " " + "high dollar projects."
# ===========================================================================
$probe2 = $d.Content
$oldSentence2 = "There is more variability in successful campaigns. This does make sense as there are simply more successful campaigns (565) than there are failed (364). "
$found2 = $probe2.Find.Execute(
    $oldSentence2,
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $partA = "There is more variability in successful campaigns. "
    $partB = "This makes sense because there are a lot more samples for successful outcomes as well as more"
    $partSpace = " "
    $partC = "high dollar projects."
    $newText2 = $partA + $partB + $partSpace + $partC

    $run2 = $d.Range($probe2.Start, $probe2.End)
    $run2.Text = $newText2
    $run2.Bold = 1
    $run2.Bold = 0

    $offset = $run2.Start

    # Split off partA.
    $rangeA = $d.Range($offset, $offset + $partA.Length)
    $rangeA.Bold = 1
    $rangeA.Bold = 0
    $offset = $offset + $partA.Length

    # Split off partB.
    $rangeB = $d.Range($offset, $offset + $partB.Length)
    $rangeB.Bold = 1
    $rangeB.Bold = 0
    $offset = $offset + $partB.Length

    # Split off the single space so it is its own run (mirrors the
    # pre-existing trailing-space run in the source document), leaving
    # partC ("high dollar projects.") as the final, separate run.
    $rangeSpace = $d.Range($offset, $offset + $partSpace.Length)
    $rangeSpace.Bold = 1
    $rangeSpace.Bold = 0
}
